# [Feat : KSW] Add Component of WeaponInfo.json
# Adds a new "Equip" (bool / 착용여부) column (M) to the WeaponInfo sheet,
# sets it to FALSE for every weapon data row, widens column M, and updates
# the active view (zoom + selection) to match the author's saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows (1:3) for the new "Equip" column -------------------------
$ws.Range("M1").Value = "Equip"
$ws.Range("M2").Value = "bool"
$ws.Range("M3").Value = "착용여부"

# Match the formatting of the neighboring header column (L) for the new
# top header cell (bold / centered "title" style).
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data rows (4:38) ------------------------------------------------------
# Every weapon row gets an "Equip" default value of FALSE.
for ($i = 4; $i -le 38; $i++) {
    $ws.Cells.Item($i, 13).Value = $false
}

# --- Column sizing ----------------------------------------------------------
$ws.Columns.Item(13).ColumnWidth = 8.29

# --- View state: zoom + selected cell --------------------------------------
$excel.ActiveWindow.Zoom = 115
$ws.Range("E5").Select() | Out-Null
